# Adds a "policy_start_date" column (with a default value of 30 "days") to
# the Sheet1 test-data table so the new "create quotation without customer
# details" test method has a policy start date to work with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The existing table occupies A1:AA9 (header row 1, 8 data rows 2-9).
# Column AB (28) is the next free column - add the new field there.
#
# Force the cells to text ("@") BEFORE writing the values so the "30" is
# stored as text (matching the rest of the sheet, which uses style index 3
# - numFmtId 49 / text - for its "code"-like columns) instead of being
# auto-typed as a number.
$ws.Range("AB1:AB9").NumberFormat = "@"

$ws.Range("AB1").Value = "policy_start_date"
$ws.Range("AB2:AB9").Value = "30"

# Move the selection/viewport to the newly added column, same as a user
# would after typing the new column in manually.
$ws.Range("AB10").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 20

# Keep the page set up for (portrait) printing now that the sheet is wider.
$ws.PageSetup.Orientation = 1
